$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure all touched cells are formatted as Text so values (especially
# the date-like strings) are stored verbatim as plain text rather than
# being auto-converted to dates/numbers by Excel, matching the source
# data which is all plain text.
$ws.Range("A2:F6").NumberFormat = "@"

# Row 2 updates
$ws.Range("B2").Value = "2024-04-16"
$ws.Range("C2").Value = "Website"
$ws.Range("D2").Value = "5000"
$ws.Range("F2").Value = "dgfy"

# Row 3 updates
$ws.Range("A3").Value = "Ram Chaudhary"
$ws.Range("B3").Value = "2024-04-05"
$ws.Range("C3").Value = "Developer"
$ws.Range("D3").Value = "1500"
$ws.Range("E3").Value = "Cash"
$ws.Range("F3").Value = "esewa bill no 2"

# Row 4 new
$ws.Range("A4").Value = "samir"
$ws.Range("B4").Value = "2024-04-01"
$ws.Range("C4").Value = "Debops"
$ws.Range("D4").Value = "400"
$ws.Range("E4").Value = "Cash"
$ws.Range("F4").Value = "dgfy"

# Row 5 new
$ws.Range("A5").Value = "admin"
$ws.Range("B5").Value = "2024-03-31"
$ws.Range("C5").Value = "Developer"
$ws.Range("D5").Value = "490"
$ws.Range("E5").Value = "Cash"
$ws.Range("F5").Value = "fg54"

# Row 6 new
$ws.Range("A6").Value = "Asmin Dhakal"
$ws.Range("B6").Value = "2024-04-17"
$ws.Range("C6").Value = "Developer"
$ws.Range("D6").Value = "5000"
$ws.Range("E6").Value = "Online"
$ws.Range("F6").Value = "Nabil 3"

# Now that the values are safely stored as text (not auto-converted into
# dates/numbers), restore the default "Normal" style on the data cells so
# the on-disk styling matches the original (unstyled) data rows.
$ws.Range("A2:F6").Style = "Normal"
